# Add three new rows (94-96) to Sheet1, replicating the pattern of row 93:
# column A holds the next sequential date (serial day number) formatted like
# the existing date column, and columns B-J repeat the same values found in
# row 93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the date cell's formatting (number format / font / border / alignment)
# from A93 onto the new date cells A94:A96 without touching their values yet.
$ws.Range("A93").Copy()
$ws.Range("A94:A96").PasteSpecial(-4122)

# Serial date values: one day after another, continuing from row 93 (45649).
$ws.Range("A94").Value = 45650
$ws.Range("A95").Value = 45651
$ws.Range("A96").Value = 45652

# Columns B through J repeat the exact same values as row 93 for each new row.
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
$values = @(
    116.4121952,
    0.00170247,
    0.008850780000000001,
    0.06933635,
    12792.90181321,
    465.80531254,
    0.24,
    1.7904431,
    485.38834923
)

foreach ($r in 94..96) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
